$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("D2").Value = 1054
$ws.Range("E2").Value = 37
$ws.Range("F2").Value = 37
$ws.Range("G2").Value = 7
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = -8
$ws.Range("K2").Value = 1907
$ws.Range("L2").Value = 1175
$ws.Range("M2").Value = 732
$ws.Range("N2").Value = 690
$ws.Range("O2").Value = 42
$ws.Range("P2").Value = 88
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = -3
$ws.Range("S2").Value = 19
$ws.Range("T2").Value = 29
$ws.Range("U2").Value = -29
$ws.Range("V2").Value = 874
$ws.Range("W2").Value = 3.53
$ws.Range("X2").Value = -0.04
$ws.Range("Y2").Value = 1.11
$ws.Range("Z2").Value = -0.02
$ws.Range("AA2").Value = 160.63
$ws.Range("AB2").Value = 552.2
$ws.Range("AC2").Value = 38
$ws.Range("AD2").Value = 56.05
$ws.Range("AE2").Value = 3444
$ws.Range("AF2").Value = 0.61
$ws.Range("AG2").Value = 16
$ws.Range("AH2").Value = 0.74
$ws.Range("AI2").Value = 41.08
$ws.Range("AJ2").Value = 20159098

# --- Row 3 ---
$ws.Range("D3").Value = 1251
$ws.Range("E3").Value = -25
$ws.Range("F3").Value = -25
$ws.Range("G3").Value = -105
$ws.Range("H3").Value = -179
$ws.Range("I3").Value = -161
$ws.Range("J3").Value = -17
$ws.Range("K3").Value = 1994
$ws.Range("L3").Value = 1323
$ws.Range("M3").Value = 671
$ws.Range("N3").Value = 632
$ws.Range("O3").Value = 39
$ws.Range("P3").Value = 91
$ws.Range("Q3").Value = -55
$ws.Range("R3").Value = -8
$ws.Range("S3").Value = 39
$ws.Range("T3").Value = 36
$ws.Range("U3").Value = -91
$ws.Range("V3").Value = 916
$ws.Range("W3").Value = -1.98
$ws.Range("X3").Value = -14.28
$ws.Range("Y3").Value = -24.41
$ws.Range("Z3").Value = -9.16
$ws.Range("AA3").Value = 197.11
$ws.Range("AB3").Value = 327.81
$ws.Range("AC3").Value = -800
$ws.Range("AD3").Value = -3.23
$ws.Range("AE3").Value = 3156
$ws.Range("AF3").Value = 0.82
$ws.Range("AG3").Value = 16
$ws.Range("AH3").Value = 0.63
$ws.Range("AI3").Value = -2.01
$ws.Range("AJ3").Value = 20159098

# --- Row 4 ---
$ws.Range("D4").Value = 1422
$ws.Range("E4").Value = 101
$ws.Range("F4").Value = 101
$ws.Range("G4").Value = 82
$ws.Range("H4").Value = 65
$ws.Range("I4").Value = 67
$ws.Range("J4").Value = -3
$ws.Range("K4").Value = 2211
$ws.Range("L4").Value = 1264
$ws.Range("M4").Value = 947
$ws.Range("N4").Value = 910
$ws.Range("O4").Value = 37
$ws.Range("P4").Value = 116
$ws.Range("Q4").Value = 56
$ws.Range("R4").Value = -47
$ws.Range("S4").Value = 142
$ws.Range("T4").Value = 93
$ws.Range("U4").Value = -37
$ws.Range("V4").Value = 837
$ws.Range("W4").Value = 7.07
$ws.Range("X4").Value = 4.55
$ws.Range("Y4").Value = 8.75
$ws.Range("Z4").Value = 3.08
$ws.Range("AA4").Value = 133.5
$ws.Range("AB4").Value = 473.22
$ws.Range("AC4").Value = 304
$ws.Range("AD4").Value = 17.3
$ws.Range("AE4").Value = 3728
$ws.Range("AF4").Value = 1.41
$ws.Range("AG4").Value = 47
$ws.Range("AH4").Value = 0.9
$ws.Range("AI4").Value = 17.11
$ws.Range("AJ4").Value = 24540000

# --- Row 5 ---
$ws.Range("D5").Value = 1401
$ws.Range("E5").Value = 75
$ws.Range("F5").Value = 75
$ws.Range("G5").Value = 101
$ws.Range("H5").Value = 81
$ws.Range("I5").Value = 78
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 2351
$ws.Range("L5").Value = 1334
$ws.Range("M5").Value = 1017
$ws.Range("N5").Value = 977
$ws.Range("O5").Value = 40
$ws.Range("P5").Value = 120
$ws.Range("Q5").Value = 80
$ws.Range("R5").Value = -276
$ws.Range("S5").Value = 23
$ws.Range("T5").Value = 216
$ws.Range("U5").Value = -136
$ws.Range("V5").Value = 872
$ws.Range("W5").Value = 5.33
$ws.Range("X5").Value = 5.8
$ws.Range("Y5").Value = 8.279999999999999
$ws.Range("Z5").Value = 3.56
$ws.Range("AA5").Value = 131.11
$ws.Range("AB5").Value = 509.84
$ws.Range("AC5").Value = 318
$ws.Range("AD5").Value = 18.4
$ws.Range("AE5").Value = 4005
$ws.Range("AF5").Value = 1.46
$ws.Range("AG5").Value = 98
$ws.Range("AI5").Value = 30.61
$ws.Range("AJ5").Value = 24540000

# --- Row 6 ---
$ws.Range("D6").Value = 1475
$ws.Range("E6").Value = 49
$ws.Range("F6").Value = 49
$ws.Range("G6").Value = 19
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 10
$ws.Range("K6").Value = 2431
$ws.Range("L6").Value = 1500
$ws.Range("M6").Value = 931
$ws.Range("N6").Value = 897
$ws.Range("P6").Value = 123
$ws.Range("Q6").Value = 114
$ws.Range("R6").Value = -161
$ws.Range("S6").Value = 40
$ws.Range("T6").Value = 109
$ws.Range("U6").Value = 5
$ws.Range("V6").Value = 1008
$ws.Range("W6").Value = 3.35
$ws.Range("X6").Value = 0.22
$ws.Range("Y6").Value = 1.05
$ws.Range("Z6").Value = 0.14
$ws.Range("AA6").Value = 161.18
$ws.Range("AB6").Value = 457.61
$ws.Range("AC6").Value = 40
$ws.Range("AD6").Value = 128.87
$ws.Range("AE6").Value = 3677
$ws.Range("AF6").Value = 1.41
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 24540000
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()

# --- Row 7 ---
$ws.Range("D7").Value = 1520
$ws.Range("E7").Value = -35
$ws.Range("G7").Value = -75
$ws.Range("H7").Value = -60
$ws.Range("I7").Value = -58
$ws.Range("K7").Value = 2564
$ws.Range("L7").Value = 1667
$ws.Range("M7").Value = 897
$ws.Range("N7").Value = 865
$ws.Range("P7").Value = 123
$ws.Range("Q7").Value = 42
$ws.Range("R7").Value = -168
$ws.Range("S7").Value = 201
$ws.Range("T7").Value = 114
$ws.Range("U7").Value = -72
$ws.Range("W7").Value = -2.3
$ws.Range("X7").Value = -3.95
$ws.Range("Y7").Value = -6.58
$ws.Range("Z7").Value = -2.4
$ws.Range("AA7").Value = 185.84
$ws.Range("AC7").Value = -236
$ws.Range("AD7").Value = -19.7
$ws.Range("AE7").Value = 3546
$ws.Range("AF7").Value = 1.31
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0
$ws.Range("AI7").ClearContents()

# --- Row 8 ---
$ws.Range("D8").Value = 1600
$ws.Range("E8").Value = 90
$ws.Range("G8").Value = 61
$ws.Range("H8").Value = 55
$ws.Range("I8").Value = 54
$ws.Range("K8").Value = 2617
$ws.Range("L8").Value = 1665
$ws.Range("M8").Value = 953
$ws.Range("N8").Value = 919
$ws.Range("P8").Value = 123
$ws.Range("Q8").Value = 86
$ws.Range("R8").Value = -4
$ws.Range("S8").Value = -28
$ws.Range("T8").Value = 0
$ws.Range("U8").Value = 86
$ws.Range("W8").Value = 5.63
$ws.Range("X8").Value = 3.44
$ws.Range("Y8").Value = 6.05
$ws.Range("Z8").Value = 2.12
$ws.Range("AA8").Value = 174.71
$ws.Range("AC8").Value = 220
$ws.Range("AD8").Value = 21.15
$ws.Range("AE8").Value = 3768
$ws.Range("AF8").Value = 1.24
$ws.Range("AG8").Value = 0
$ws.Range("AH8").Value = 0
$ws.Range("AI8").ClearContents()

# --- Row 9 ---
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
